$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the component name in A6 to include the full part number.
$ws.Range("A6").Value = "STM32 Nucleo64 STM32L476RGT6"

# Column A was bestFit previously; re-fit it now that the text is longer.
[void]$ws.Range("A1").EntireColumn.AutoFit()

# Move the active selection to B10, matching the saved view state.
[void]$ws.Range("B10").Select()
